# Commit: "Remove interpolation row from some of the moorings data and
# switch some of the structures to LineTables."
#
# For dds/constants.xlsx the LineTable portion of that commit shows up as:
#   - ROOT sheet, column B ("Structure" type), rows 9-14 (the six
#     constants.* structures that used to be TableDataColumn) are switched
#     to LineTableColumn.
#   - The ROOT sheet becomes the active/selected tab (it was "Units"
#     before), with the active cell now on B14.

$wb = $excel.ActiveWorkbook

$root = $wb.Worksheets.Item("ROOT")

# Switch the six TableDataColumn structures to LineTableColumn.
$root.Range("B9:B14").Value = "LineTableColumn"

# Make ROOT the active sheet/tab, with B14 as the active selection -
# this mirrors the workbook view moving off "Units" (which loses
# tabSelected) and onto "ROOT".
$root.Select()
$root.Range("B14").Select()
